$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "61.136.35"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = "  -0.69%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "2.389.31"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = "  -3.59%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "548.87"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = "  -0.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "142.11"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = "  -2.63%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.541"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = "  -10.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "2.386.44"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = "  -3.66%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "0.106"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = "  -1.58%  "
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = "  +0.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "5.25"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = "  -3.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "0.349"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = "  -2.32%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "25.50"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = "  -2.61%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "2.822.13"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = "  -3.50%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = "  -0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "61.014.25"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = "  -0.71%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = "2.388.70"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = "  -3.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "10.73"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = "  -3.97%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "4.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = "  -1.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "319.70"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = "  -0.69%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "6.79"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = "  -2.47%  "
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value2 = "SuiNetwork"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value2 = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "1.96"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = "  +7.85%  "
$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value2 = "Dai"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value2 = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "1.00"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = "  +0.09%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "63.71"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = "  -0.52%  "
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value2 = "Aptos"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value2 = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "8.15"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = "  +6.12%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.0₃0948"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = "  -4.53%  "
$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value2 = "Binance-PegBSC-USD"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.00"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = "  +0.08%  "
$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value2 = "WrappedeETH"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.502.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = "  -3.58%  "
$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value2 = "Bittensor"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value2 = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "533.68"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = "  +0.40%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.44"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = "  -4.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "8.13"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = "  -3.51%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "0.146"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = "  -3.78%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.84"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = "  -3.99%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "1.60"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = "  +1.11%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = "  +0.00%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "5.60"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = "  -5.54%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "4.73"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = "  -3.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.378"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = "  -1.39%  "
$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value2 = "Stacks"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "1.83"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = "  +5.17%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value2 = "EthereumClassic"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "18.14"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = "  -2.21%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = "  +0.12%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "137.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = "  -7.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "40.42"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = "  +0.08%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "2.29"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = "  -2.46%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "142.63"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = "  -3.65%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "3.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "20.34"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = "  -2.25%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = "0.0521"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = "  -2.48%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.578"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = "  -3.16%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.0909"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = "  -4.17%  "
